$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new key/value rows (total_errors / total_score) ---
# They land right after the existing "score" row (row 20) and before
# "health_warning" (old row 21), becoming new rows 21 and 22.
$ws.Range("A21:A22").EntireRow.Insert() | Out-Null

$ws.Cells.Item(21, 1).Value = "total_errors"
$ws.Cells.Item(21, 2).Value = "TOTAL ERROR"

$ws.Cells.Item(22, 1).Value = "total_score"
$ws.Cells.Item(22, 2).Value = "TOTAL SCORE"

# --- Update existing text values ---

# attack_instruct_success_0 (row 66, after the insert above)
$ws.Cells.Item(66, 2).Value = "Excellent! We’ve managed to clear the blobs!"

# level1_intro_1 (row 60, after the insert above)
$ws.Cells.Item(60, 2).Value = "In order to beat the mega blob, we must merge all the blobs into one final quotient blob."

# end_desc VoiceDuration (column C) changes from 7 to 5 (row 93, after the insert above)
$ws.Cells.Item(93, 3).Value = 5

# --- Update the view selection to reflect the edited area ---
$ws.Range("C93").Select() | Out-Null
